# KKBOX autoint corr map + ieee_fraud_small_data_test
# - Rename "Sheet2" to "model reference"
# - Update Sheet1's selection/scroll position (drop the old topLeftCell=A25
#   scroll, select J8)
# - Update "model reference" sheet's selection (B25 -> B21)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Rename Sheet2 -> "model reference"
$ws2.Name = "model reference"

# Update the selection on the "model reference" sheet first (it must not
# end up as the active/tab-selected sheet).
$ws2.Activate()
$ws2.Range("B21").Select()

# Finish on Sheet1 so it remains the active tab, with the new selection.
$ws1.Activate()
$ws1.Range("J8").Select()
